$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per the scraped-site refresh.
# Cells whose new text looks like a plain number (e.g. "1.00", "8.00")
# are entered with a leading apostrophe so Excel keeps them as literal
# text (matching the workbook's existing inline-string price format)
# instead of silently converting them to numeric values.

$ws.Range('D2').Value = '62.105.21'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '3.411.48'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''575.91'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').Value = '''147.62'
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('D9').Value = '''8.00'
$ws.Range('E9').Value = '  +4.43%  '
$ws.Range('D10').Value = '''0.123'
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').Value = '''0.413'
$ws.Range('E11').Value = '  +3.18%  '
$ws.Range('D12').Value = '4.008.14'
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('D14').Value = '''28.25'
$ws.Range('E14').Value = '  -4.69%  '
$ws.Range('D15').Value = '3.424.14'
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').Value = '62.316.62'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').Value = '''6.35'
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('D19').Value = '''14.38'
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('D20').Value = '''8.90'
$ws.Range('E20').Value = '  -2.82%  '
$ws.Range('D21').Value = '''381.27'
$ws.Range('E21').Value = '  -1.92%  '
$ws.Range('D22').Value = '''0.565'
$ws.Range('E22').Value = '  +1.72%  '
$ws.Range('D23').Value = '''74.69'
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = '3.585.74'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('D26').Value = '''0.0000111'
$ws.Range('E26').Value = '  -3.38%  '
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('D28').Value = '''7.61'
$ws.Range('E28').Value = '  +1.57%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '''7.91'
$ws.Range('E30').Value = '  -2.64%  '
$ws.Range('D31').Value = '''2.12'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = '''1.00'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('E33').Value = '  -3.11%  '
$ws.Range('D34').Value = '''23.10'
$ws.Range('E34').Value = '  -2.20%  '
$ws.Range('D35').Value = '''5.45'
$ws.Range('E35').Value = '  +3.86%  '
$ws.Range('D36').Value = '''1.62'
$ws.Range('E36').Value = '  +3.48%  '
$ws.Range('D37').Value = '''31.37'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = '''168.91'
$ws.Range('E38').Value = '  -0.41%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = '''6.87'
$ws.Range('E39').Value = '  -2.35%  '
$ws.Range('D40').Value = '3.457.13'
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = '''0.0782'
$ws.Range('E41').Value = '  +3.88%  '
$ws.Range('D42').Value = '''0.783'
$ws.Range('E42').Value = '  -2.00%  '
$ws.Range('D43').Value = '''42.41'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').Value = '''4.36'
$ws.Range('E44').Value = '  -2.21%  '
$ws.Range('D45').Value = '''1.67'
$ws.Range('E45').Value = '  -2.21%  '
$ws.Range('E46').Value = '  -3.33%  '
$ws.Range('D47').Value = '2.539.45'
$ws.Range('E47').Value = '  -2.29%  '
$ws.Range('D48').Value = '''6.90'
$ws.Range('E48').Value = '  +2.98%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').Value = '''1.00'
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = '''2.19'
$ws.Range('E50').Value = '  -2.65%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '''22.44'
$ws.Range('E51').Value = '  -1.83%  '
